$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The update swaps the full data payload (columns B through AD) between
# certain pairs of rows, while leaving column A (the sequential row id)
# untouched. Swap each listed pair using a temporary array buffer.

$pairs = @(
    @(106,107),
    @(137,138),
    @(139,140)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value()
    $vals2 = $rng2.Value()

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}
